$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reindex the output to match desired order for cut and paste into excel
# Row 2 (Checked)
$ws.Range("B2").Value = 0.0880503144654088
$ws.Range("C2").Value = 0.1257861635220126
$ws.Range("D2").Value = 0.1132075471698113
$ws.Range("E2").Value = 0.2327044025157233

# Row 3 (Unchecked)
$ws.Range("B3").Value = 0.1320754716981132
$ws.Range("C3").Value = 0.1320754716981132
$ws.Range("D3").Value = 0.0440251572327044
$ws.Range("E3").Value = 0.1320754716981132
